$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B so it is (as close as the COM ColumnWidth API allows to)
# the same width as columns A and C. ColumnWidth is stored/rounded by Excel
# to whole pixels (Maximum Digit Width = 7px for Calibri 11), so 10.8 is the
# input that lands on the pixel bucket nearest to the 11.7109375 target.
$ws.Columns.Item(2).ColumnWidth = 10.8

# Update the three data values for this case
$ws.Range("A1").Value = 158.02969566758335
$ws.Range("B1").Value = 3.4556599640933605
$ws.Range("C1").Value = 5.4297129810828437
